$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.445.75"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.667.95"
$ws.Range("E3").Value = "  -2.71%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.96"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.51"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.667.20"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.78"
$ws.Range("E14").Value = "  -3.85%  "
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("E16").Value = "  -3.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.387.19"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.669.02"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.81"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.03"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("E22").Value = "  -4.11%  "
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("E24").Value = "  -5.19%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.72"
$ws.Range("E26").Value = "  -4.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.811.22"
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("E29").Value = "  -4.18%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "548.51"
$ws.Range("E31").Value = "  -8.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  -4.30%  "
$ws.Range("E33").Value = "  -4.89%  "
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  -5.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.43"
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("E39").Value = "  -4.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.371"
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("E41").Value = "  -4.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.24"
$ws.Range("E42").Value = "  -4.70%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  -7.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.30"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.588"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.37"
$ws.Range("E49").Value = "  -3.41%  "
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("E51").Value = "  -4.71%  "
